$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Airline1")
$ws.Range("C14").Value = " "
$ws.Range("I3").Select()
